$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row above the "Total hours Spent" merged block -----------
# Original sheet had rows ...37 (data), 39 (Total, merged 39:41) — row 38 was
# a genuine gap (no cells). Inserting a row at 39 pushes the Total block down
# to 40:42 while leaving the (still-empty) row 38 untouched, ready to be
# filled in below. This matches the target layout exactly.
$ws.Rows("39:39").Insert() | Out-Null

# --- Fill in the new timeline entry on (now real) row 38 -------------------
# Copy formatting from matching columns on existing rows so the new cells
# pick up the same cell styles (centered text, date format, etc.) rather
# than defaulting to "Normal".
$ws.Range("A37").Copy()
$ws.Range("A38").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B33").Copy()                # B33 carries the date-format style
$ws.Range("B38").PasteSpecial(-4122)

$ws.Range("C37").Copy()
$ws.Range("C38").PasteSpecial(-4122)

$ws.Range("D37").Copy()
$ws.Range("D38").PasteSpecial(-4122)

$excel.CutCopyMode = 0

$ws.Range("A38").Value2 = 30
$ws.Range("B38").Value2 = 45570
$ws.Range("C38").Value2 = 5
$ws.Range("D38").Value2 = "Changed User Protected Details update with process phases for better FE controll"

# --- Update the saved view/selection state ----------------------------------
$ws.Range("C25").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1

Write-Output "done"
